# Generate Report for Handoff
# Updates the Priority and Latest Handoff Datetime columns for the four
# "Ready for handoff" rows (26621909-..., c228934c-..., c799f5b3-..., e768144c-...)
# on both the "zh-cn" and "de-de" localization-status sheets.

$wb = $excel.ActiveWorkbook

$rows = @(4, 5, 6, 7)

# zh-cn sheet: priority bumped to "ht" and handoff timestamp refreshed
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-19 18:40:20"
}

# de-de sheet: same priority bump, with its own refreshed handoff timestamp
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-19 18:40:24"
}

# Overview sheet: "Latest HO Xliff Generate Date" reflects the newest of the
# per-language handoff timestamps (de-de's 18:40:24 is the latest of the two)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-19 18:40:24"
}
